# Insert a new row above row 193; this shifts the existing row 193 (and all
# rows below it) down by one, growing the used range from A1:T261 to A1:T262.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new weekly price observation.
# Columns A, B, C, E, F, G, H, I, J, K, R are identical to the record that
# used to occupy row 193 (now shifted to row 194), so copy them across;
# columns D, L, M, N, O, P, Q, S, T carry new values for this record.
$ws.Range("A193").Value = 11
$ws.Range("B193").Value = "Vega Monumental Concepción"
$ws.Range("C193").Value = "Bíobío"
$ws.Range("D193").Value = 45027
$ws.Range("E193").Value = 8
$ws.Range("F193").Value = "Fruta"
$ws.Range("G193").Value = 100108
$ws.Range("H193").Value = "Tropicales y subtropicales"
$ws.Range("I193").Value = 100108005
$ws.Range("J193").Value = "Piña"
$ws.Range("K193").Value = "Caramelo"
$ws.Range("L193").Value = "Primera"
$ws.Range("M193").Value = 150
$ws.Range("N193").Value = 19000
$ws.Range("O193").Value = 20000
$ws.Range("P193").Value = 19667
$ws.Range("Q193").Value = "$/caja 12 unidades"
$ws.Range("R193").Value = "Ecuador"
$ws.Range("S193").Value = 1639
$ws.Range("T193").Value = 12
